$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '42.606.02'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +2.19%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.296.31'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +1.38%  '
$ws.Range("E4").Value = '  +0.07%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '307.97'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +1.28%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '97.06'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +5.52%  '
$ws.Range("E7").Value = '  +0.70%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  +3.72%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '36.31'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +12.78%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0806'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +1.21%  '
$ws.Range("E12").Value = '  -1.53%  '
$ws.Range("E13").Value = '  +2.34%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '2.651.90'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +1.41%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '14.65'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +2.98%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '2.301.78'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +2.42%  '
$ws.Range("E17").Value = '  +5.42%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '42.518.36'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +2.18%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '12.75'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +2.06%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '0.0₃0922'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +2.04%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '6.03'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +2.03%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '67.95'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +1.82%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '243.33'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +1.17%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '2.62'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.67%  '
$ws.Range("E25").Value = '  +2.19%  '
$ws.Range("E26").Value = '  -0.13%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '24.07'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +0.33%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '36.84'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +7.92%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '9.59'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.94%  '
$ws.Range("E30").Value = '  +2.26%  '
$ws.Range("E31").Value = '  +0.65%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '5.32'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +1.89%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +0.11%  '
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '3.13'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +4.16%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.0753'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +1.26%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '17.36'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +4.33%  '
$ws.Range("E37").Value = '  +3.37%  '
$ws.Range("E38").Value = '  +5.29%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.38'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +0.55%  '
$ws.Range("E40").Value = '  -0.21%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '4.20'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +6.74%  '
$ws.Range("E42").Value = '  +19.77%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '2.008.76'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -2.01%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '19.41'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -0.61%  '
$ws.Range("E45").Value = '  +3.14%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '10.30'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -1.40%  '
$ws.Range("E47").Value = '  +5.23%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '53.99'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +4.32%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.55'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +0.95%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '72.76'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -0.03%  '
$ws.Range("E51").Value = '  -0.39%  '
